$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Value to place in column B for rows 2-8 (date/time serial 37399.552777777775)
$dateValue = 37399.552777777775

# Update existing row 2, and add rows 3-8 with the same data as row 2 (A=3635260, B=date)
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = 3635260
    $ws.Cells.Item($r, 2).Value = $dateValue
}

# Update selection to B4 as per the diff
$ws.Range("B4").Select()
